$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.170.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.06%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.782.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.91%  '
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.83%  '
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3934'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.40%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3419'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.97%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.86'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.36%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.190'
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07439'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.74%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.12%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.60'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.07%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.437'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.08%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.780.32'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.91%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.084'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.20%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001091'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.50%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06667'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.10%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.28'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.22%  '
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.00%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.64'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.48%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.479'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.13%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.180.78'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.01%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.34'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.75%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.379'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.78%  '
# Row 26
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.25'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.50%  '
# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.502'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -7.61%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.452'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.50%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '156.23'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.28%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.980.23'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.97%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.60'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.35%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.967'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.89%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.981'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.97%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08726'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.21%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -7.30%  '
# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.80%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02385'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.88%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.396'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.92%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6791'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.80%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06382'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.54%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2203'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.83%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.238'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.80%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.442'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.86%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.25'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.17%  '
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.05%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6387'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.98%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.856'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.85%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.134'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.42%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.95'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.33%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07097'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.11%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.50'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.42%  '
